# Daily attendance processing - 2025-12-20 10:54:26
# Rotate the "Recorded By" (column G) comma-separated list left by one
# position (move the first entry to the end) for every data row.
# Cells with only a single recorder entry are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ','
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }

        if ($parts.Length -gt 1) {
            $rotated = $parts[1..($parts.Length - 1)] + $parts[0]
            $newVal = [string]::Join(", ", $rotated)
            $cell.Value = $newVal
        }
    }
}
